# Applies the LOQ4061 syllabus update:
#  - "Ativação:" date changed from 01/01/2020 to 01/01/2022
#  - Removal of the "4)" clause (and related trailing text) from the
#    "Programa resumido:", "Short syllabus:", "Programa:" and "Syllabus:"
#    entries (both the Portuguese column B and the English column C copies
#    hold the same text in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOQ4061")

# --- Ativação: 01/01/2020 -> 01/01/2022 ---
# A plain  $range.Value = "01/01/2022"  assignment would be auto-recognised
# by Excel as a date literal (since the cell's number format is General) and
# get silently converted to a date serial number instead of staying as the
# literal text "01/01/2022". To keep it as plain text (matching the original
# file, which stores it as a shared string) without touching the cell's
# number format/style, we compute the literal text in a scratch cell via a
# formula (whose *result* is never re-interpreted as a date), copy just the
# computed value over with Paste Special "Values", then wipe the scratch
# cell again.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="01/01/2022"'
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# --- Programa resumido: ---
$resumidoNew = "Perfis de temperaturas em barras de seção circular; 2) Transferência de calor por convecção; 3) Determinação do coeficiente de difusão em sistemas gás-líquido;"
$ws.Range("B14").Value = $resumidoNew
$ws.Range("C14").Value = $resumidoNew

# --- Short syllabus: ---
$shortSyllabusNew = "1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems;"
$ws.Range("B15").Value = $shortSyllabusNew
$ws.Range("C15").Value = $shortSyllabusNew

# --- Programa: ---
$programaNew = "1) Perfis de temperaturas em barras de seção circular: processos envolvendo condução e convecção em barras de vários materiais e diferentes dimensões. Aplicação do princípio das aletas; 2) Transferência de calor por convecção: medidas da variação de temperatura em corpos de várias geometrias e materiais diferentes e comparação com a análise concentrada para regime transiente; 3) Determinação do coeficiente de difusão em sistemas gás-líquido: avaliação da transferência de massa entre ar e líquidos empregando tubos horizontais (célula de Stefan) em regime pseudo-estacionário;"
$ws.Range("B16").Value = $programaNew
$ws.Range("C16").Value = $programaNew

# --- Syllabus: ---
$syllabusNew = "1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state;"
$ws.Range("B17").Value = $syllabusNew
$ws.Range("C17").Value = $syllabusNew

$wb.Save()
